{"js": "const body = context.document.body;\n\n// ------------------------------------------------------------------\n// 1) \"So, lets start with undefined...\" -> \"So, Let's start with undefined...\"\n//    (curly apostrophe, capital L)\n// ------------------------------------------------------------------\nconst introResults = body.search(\"lets start with undefined\", { matchCase: false });\nintroResults.load(\"text\");\nawait context.sync();\nintroResults.items[0].insertText(\"Let\\u2019s start with undefined\", Word.InsertLocation.replace);\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 2) Italicize the three code-snippet paragraphs + the \"Undefined\" list item\n// ------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst italicTargets = [\"var name;\", \"var age;\", \"age;\", \"Undefined \"];\nfor (const p of paragraphs.items) {\n  if (italicTargets.includes(p.text)) {\n    p.font.italic = true;\n  }\n}\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 3) Rewrite the \"but we do not set it...\" sentence\n// ------------------------------------------------------------------\nconst sentenceResults = body.search(\n  \"but we do not set it to any value or strings, its declared but not initialized to a value, then its considered undefined. \",\n  { matchCase: true }\n);\nsentenceResults.load(\"text\");\nawait context.sync();\nsentenceResults.items[0].insertText(\n  \"but we do not set it to any value or strings, then its declared but not initialized to a value, then it is considered undefined.  \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 4) Move the \"_GoBack\" bookmark into the \"Null is another way of...\" paragraph\n//    (right after \"Null is another way of e\")\n// ------------------------------------------------------------------\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n} catch (e) {\n  // no-op if it doesn't currently exist\n}\nawait context.sync();\n\nconst bookmarkAnchor = body.search(\"Null is another way of e\", { matchCase: true });\nbookmarkAnchor.load(\"text\");\nawait context.sync();\nconst afterAnchor = bookmarkAnchor.items[0].getRange(\"End\");\nafterAnchor.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop script applying the \"Null & Undefined\" edits.\n$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) \"So, lets start with undefined...\" -> \"So, Let's start with undefined...\"\n# ------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\"lets start with undefined\", $false, $false, $false, $false, $false, $true, 1, $false, \"Let's start with undefined\", 2) | Out-Null\n\n# ------------------------------------------------------------------\n# 2) Italicize the three code-snippet paragraphs + the \"Undefined\" list item\n# ------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"var name;\" -or $t -eq \"var age;\" -or $t -eq \"age;\" -or $t -eq \"Undefined \") {\n        $p.Range.Font.Italic = 1\n        $p.Range.Font.ItalicBi = 1\n    }\n}\n\n# ------------------------------------------------------------------\n# 3) Rewrite the \"but we do not set it...\" sentence\n# ------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\"but we do not set it to any value or strings, its declared but not initialized to a value, then its considered undefined. \", $false, $false, $false, $false, $false, $true, 1, $false, \"but we do not set it to any value or strings, then its declared but not initialized to a value, then it is considered undefined.  \", 2) | Out-Null\n\n# ------------------------------------------------------------------\n# 4) Move the \"_GoBack\" bookmark into the \"Null is another way of...\" paragraph\n#    (right after \"Null is another way of e\")\n# ------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\"Null is another way of e\") | Out-Null\n$insertPoint = $d.Range($rng.End, $rng.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertPoint) | Out-Null\n\nWrite-Output \"done\"\n"}
